$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "+" marker to D4 (русско таджикский словарь) - now used in the WordDialog
$ws.Range("D4").Value = "+"

# Remove EN-TJ / TJ-EN words from the DB (rows 7, 9 and 13)
# Row 7: "таджикско английский переводчик"
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = ""

# Row 9: "переводчик с английского на таджикский"
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = ""

# Row 13: "таджикско английский словарь"
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("D13").Value = ""

# Update the active selection (WordDialog cell) shown when the file was saved
$ws.Range("J11").Select()
